$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lines")

# Make "Lines" the active/selected sheet (mirrors tabSelected moving from
# "Picks" to "Lines" and workbookView activeTab="1").
$ws.Activate()

# New header cells for columns C ("League") and D ("Year").
$ws.Range("C1").Value = "League"
$ws.Range("D1").Value = "Year"

# Fill League/Year values for the 30 existing data rows (rows 2-31).
$ws.Range("C2:C31").Value = "MLB"
$ws.Range("D2:D31").Value = 2023

# Match the formatting used by the existing "League"-style header/label
# cells elsewhere in the workbook (black Arial 10) on the new C/D header
# cells and the new League (C) column.
$hdr = $ws.Range("C1:D1")
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 10
$hdr.Font.Color = 0

$league = $ws.Range("C2:C31")
$league.Font.Name = "Arial"
$league.Font.Size = 10
$league.Font.Color = 0

# Update the sheet's remembered selection to match the new active cell.
$ws.Range("E24").Select() | Out-Null
